$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from column N into new column O (row 3, and rows 5-25) ---
# Row 3 header cell (skip row 4, which stays untouched / has no column O data)
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

# Rows 5 through 25 (data rows), formats copied in one shot, values set after
$ws.Range("N5:N25").Copy()
$ws.Range("O5:O25").PasteSpecial(-4122)

# --- Set the new column O values (2021 figures) ---
$ws.Range("O3").Value = 2021

$ws.Range("O5").Value = 2148.1999999999998
$ws.Range("O6").Value = 109.5
$ws.Range("O7").Value = 210.1
$ws.Range("O8").Value = 196
$ws.Range("O9").Value = 209
$ws.Range("O10").Value = 300.2
$ws.Range("O11").Value = 302.89999999999998
$ws.Range("O12").Value = 786
$ws.Range("O13").Value = 27.7
$ws.Range("O14").Value = 6.8
# O15 intentionally stays blank (just formatted, like N15)
$ws.Range("O16").Value = 26.9
$ws.Range("O17").Value = 15.9
$ws.Range("O18").Value = 21.7
$ws.Range("O19").Value = 29.9
$ws.Range("O20").Value = 30.2
$ws.Range("O21").Value = 24
$ws.Range("O22").Value = 31.6
$ws.Range("O23").Value = 30.3
$ws.Range("O24").Value = 20.7
$ws.Range("O25").Value = 12

# --- Update the selected / active cell shown when the sheet is opened ---
$ws.Range("Q20").Select()
